$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.987.77"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.48"
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.01"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3893"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3849"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.12"
$ws.Range("E9").Value = "  +4.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.359"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08483"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.94"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.210"
$ws.Range("E14").Value = "  +3.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.008"
$ws.Range("E15").Value = "  +8.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001312"
$ws.Range("E16").Value = "  +3.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.660.05"
$ws.Range("E17").Value = "  +2.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.36"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06991"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.88"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.953"
$ws.Range("E21").Value = "  +2.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.66"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.976.91"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.487"
$ws.Range("E25").Value = "  +3.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.071"
$ws.Range("E26").Value = "  +11.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.12"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.42"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("B29").Value = "HuobiToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.361"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "139.65"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.841"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.486"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.837.65"
$ws.Range("E33").Value = "  +2.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.044"
$ws.Range("E34").Value = "  +8.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08118"
$ws.Range("E35").Value = "  +1.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02991"
$ws.Range("E36").Value = "  +4.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.06"
$ws.Range("E37").Value = "  +6.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.715"
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2698"
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09139"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.61"
$ws.Range("E41").Value = "  +2.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7559"
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.420"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.49"
$ws.Range("E44").Value = "  +4.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7004"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.481"
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.081"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08264"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.93"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.238"
$ws.Range("E51").Value = "  -0.33%  "

Write-Host "Applied 101 cell updates"
